$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Experimental" row (row 7) now carries the literal text value "true"
# (leading apostrophe forces text instead of a boolean literal)
$ws.Range("B7").Value = "'true"

# Re-apply the plain (non quote-prefixed) cell format so B7 keeps the
# same visual style as the rest of the table after the text-forcing entry
$ws.Range("B8").Copy()
$ws.Range("B7").PasteSpecial(-4122)  # xlPasteFormats

# "Date" row (row 8) gets refreshed to the new generation timestamp
$ws.Range("B8").Value = "2025-01-28T15:58:19+00:00"
